$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.150.60"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.584.45"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.64"
$ws.Range("E5").Value = "  +2.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.62"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.599"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.593.20"
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.68"
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  +3.08%  "
$ws.Range("E12").Value = "  +8.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.342"
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.034.07"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.198.13"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.47"
$ws.Range("E16").Value = "  +8.16%  "
$ws.Range("E17").Value = "  +4.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.581.15"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.50"
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "336.84"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.18"
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.19"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.06"
$ws.Range("E24").Value = "  -3.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.449"
$ws.Range("E25").Value = "  +5.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +2.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.24"
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0782"
$ws.Range("E29").Value = "  +3.99%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.67"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.57"
$ws.Range("E32").Value = "  +2.88%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.04"
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.00"
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.03"
$ws.Range("E35").Value = "  +3.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.884"
$ws.Range("E36").Value = "  +8.20%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.13"
$ws.Range("E37").Value = "  +2.23%  "
$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.873"
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.50"
$ws.Range("E39").Value = "  +2.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.82"
$ws.Range("E40").Value = "  -1.03%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.65"
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "291.63"
$ws.Range("E42").Value = "  +4.22%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0971"
$ws.Range("E44").Value = "  +2.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.598"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0535"
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.62"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.14"
$ws.Range("E48").Value = "  +2.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.14"
$ws.Range("E49").Value = "  +8.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0232"
$ws.Range("E50").Value = "  +2.76%  "
$ws.Range("E51").Value = "  +4.40%  "
